$d = $word.ActiveDocument
$sec = $d.Sections.Item(1)

# Footer 1: Pearson Edexcel logo (docPr id="3") -> rename image2.png to image1.png
$f1 = $sec.Footers.Item(1)
if ($f1.Exists -and $f1.Range.InlineShapes.Count -ge 1) {
    $pearson1 = $f1.Range.InlineShapes.Item(1)
    if ($pearson1.AlternativeText -eq "Y:\Together Design\Pearson Edexcel PowerPoint amends\Assets\PearsonLogo.png") {
        $pearson1.Name = "image1.png"
        Write-Output "Renamed footer 1 Pearson logo to image1.png"
    }
}

# Footer 2: Pearson Edexcel logo (docPr id="2") -> rename image2.png to image1.png
$f2 = $sec.Footers.Item(2)
if ($f2.Exists -and $f2.Range.InlineShapes.Count -ge 1) {
    $pearson2 = $f2.Range.InlineShapes.Item(1)
    if ($pearson2.AlternativeText -eq "Y:\Together Design\Pearson Edexcel PowerPoint amends\Assets\PearsonLogo.png") {
        $pearson2.Name = "image1.png"
        Write-Output "Renamed footer 2 Pearson logo to image1.png"
    }
}

# Header(s): BTec logo (docPr id="1") -> rename image1.jpg to image2.jpg
for ($hi = 1; $hi -le $sec.Headers.Count; $hi++) {
    $h = $sec.Headers.Item($hi)
    if ($h.Exists -and $h.Range.InlineShapes.Count -ge 1) {
        $btec = $h.Range.InlineShapes.Item(1)
        if ($btec.AlternativeText -eq "BTec_Logo-Orange") {
            $btec.Name = "image2.jpg"
            Write-Output "Renamed header $hi BTec logo to image2.jpg"
        }
    }
}
